# Generate Report for Handoff
# Updates the localization-status report after a fresh handoff xliff
# generation run: bumps the "Latest HO Xliff Generate Date" / "Latest
# Handoff Datetime" timestamps for the affected files and marks their
# Priority as "ht" (handoff type) on the per-language sheets.

$wb = $excel.ActiveWorkbook

$rows = @(7, 8, 9, 10, 11, 13)

# Overview sheet: bump the "Latest HO Xliff Generate Date" column (G)
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Cells.Item($r, 7).Value = "2016-08-24 06:21:55"
}

# zh-cn sheet: bump "Latest Handoff Datetime" (H) and set Priority (E) to "ht"
$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZhCn.Cells.Item($r, 5).Value = "ht"
    $wsZhCn.Cells.Item($r, 8).Value = "2016-08-24 06:21:49"
}

# de-de sheet: set Priority (E) to "ht". Its "Latest Handoff Datetime" (H)
# for these rows happens to already equal the Overview's old timestamp, so
# it moves in lockstep to the same new value used above.
$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDeDe.Cells.Item($r, 5).Value = "ht"
    $wsDeDe.Cells.Item($r, 8).Value = "2016-08-24 06:21:55"
}
